$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Status text update: "Ready for handoff" -> "Handed back: in sync with en-US"
# Status (column C) is shared across the zh-cn and de-de sheets.
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")
$ov = $wb.Worksheets.Item("Overview")

$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Helper values
# ---------------------------------------------------------------------------
$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce15f614f1085b4ba376681c5a9d6604108d870b/e2e/a.md"
$bUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce15f614f1085b4ba376681c5a9d6604108d870b/e2e/b.md"
$linkColor = 15570276  # BGR for RGB FF6495ED, matching the workbook's HyperLink style

# ---------------------------------------------------------------------------
# zh-cn sheet: fill in "Latest Target File" (I), "Latest Handback File" (J)
# and "Latest Handback DateTime" (K) for both data rows.
# ---------------------------------------------------------------------------
$zh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-03 04:41:00"
$zh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-03 04:41:00"

# Add hyperlinks for the newly filled "Latest Target File" cells without
# touching the existing A2/A3 hyperlinks (leave their style/relationship
# untouched).
$zh.Range("I2").Value = "a.md"
$zh.Range("I3").Value = "a.md"
$zh.Hyperlinks.Add($zh.Range("I2"), $aUrl, "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("I3"), $aUrl, "", "", "a.md")
$zh.Range("I2").Font.Underline = 2
$zh.Range("I2").Font.Color = $linkColor
$zh.Range("I3").Font.Underline = 2
$zh.Range("I3").Font.Color = $linkColor

# ---------------------------------------------------------------------------
# de-de sheet: same fields, with the de-de handback file/date.
# ---------------------------------------------------------------------------
$de.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$de.Range("K2").Value = "2016-09-03 04:41:12"
$de.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$de.Range("K3").Value = "2016-09-03 04:41:12"

$de.Range("I2").Value = "a.md"
$de.Range("I3").Value = "a.md"
$de.Hyperlinks.Add($de.Range("I2"), $aUrl, "", "", "a.md")
$de.Hyperlinks.Add($de.Range("I3"), $aUrl, "", "", "a.md")
$de.Range("I2").Font.Underline = 2
$de.Range("I2").Font.Color = $linkColor
$de.Range("I3").Font.Underline = 2
$de.Range("I3").Font.Color = $linkColor

# ---------------------------------------------------------------------------
# Column widths: widen the columns that now hold longer text.
# ---------------------------------------------------------------------------
$ov.Columns.Item(5).ColumnWidth = 29.16   # zh-cn column
$ov.Columns.Item(6).ColumnWidth = 29.16   # de-de column

$zh.Columns.Item(3).ColumnWidth = 29.16   # Status
$zh.Columns.Item(10).ColumnWidth = 39.1667 # Latest Handback File

$de.Columns.Item(3).ColumnWidth = 29.16   # Status
$de.Columns.Item(10).ColumnWidth = 39.1667 # Latest Handback File
